$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 211; existing rows 211:290 shift down to 212:291
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new weekly price entry
$ws.Range("A211").Value = 5
$ws.Range("B211").Value = "Macroferia Regional de Talca"
$ws.Range("C211").Value = "Maule"
$ws.Range("D211").Value = 44825
$ws.Range("D211").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E211").Value = 7
$ws.Range("F211").Value = "Fruta"
$ws.Range("G211").Value = 100108
$ws.Range("H211").Value = "Tropicales y subtropicales"
$ws.Range("I211").Value = 100108005
$ws.Range("J211").Value = "Piña"
$ws.Range("K211").Value = "Caramelo"
$ws.Range("L211").Value = "Tercera"
$ws.Range("M211").Value = 320
$ws.Range("N211").Value = 21000
$ws.Range("O211").Value = 21000
$ws.Range("P211").Value = 21000
$ws.Range("Q211").Value = "$/caja 16 unidades"
$ws.Range("R211").Value = "Ecuador"
$ws.Range("S211").Value = 1312
$ws.Range("T211").Value = 16
